$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = -7.329199999999991
$ws.Range("A10").Value = -20.44649999999997
$ws.Range("A12").Value = -22.54920000000002
$ws.Range("C13").Value = -12.75289999999999
$ws.Range("A18").Value = -22.50890000000002
$ws.Range("D20").Value = -8.575699999999994
